# Applies the "erase fish content" edit to the CareTeam StructureDefinition
# workbook: updates the generation Date on the Metadata sheet, and on the
# Elements sheet replaces the free-text description of
# CareTeam.participant.extension with a concise "Extension" / "An Extension"
# pair, then inserts a brand new row describing the
# CareTeam.participant.extension:comment slice.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet: bump the generation timestamp
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-07-18T14:36:40+02:00"

# ---------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Row 19 (CareTeam.participant.extension) loses its long boilerplate
# description in favor of short "Extension" / "An Extension" values, and a
# couple of other supporting cells are simplified/cleared.
$ws.Range("D19").Value = ""
$ws.Range("L19").Value = "Extension"
$ws.Range("M19").Value = "An Extension"
$ws.Range("N19").ClearContents()
$ws.Range("AB19").Value = "value:url}`n"
$ws.Range("AC19").ClearContents()
$ws.Range("AE19").Value = "open"
$ws.Range("AM19").Value = ""

# Insert a brand-new row 20 for the CareTeam.participant.extension:comment
# slice. Copy formatting down from row 19 first so the new row picks up the
# same cell style used throughout the table.
$ws.Rows.Item(20).Insert()
$ws.Range("A19:AM19").Copy()
$ws.Range("A20:AM20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A20").Value = "CareTeam.participant.extension:comment"
$ws.Range("B20").Value = "CareTeam.participant.extension"
$ws.Range("C20").Value = "comment"
$ws.Range("D20").Value = ""
$ws.Range("E20").ClearContents()
$ws.Range("F20").Value = "0"
$ws.Range("G20").Value = "1"
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = "Extension {http://example.org/fhir/fish/StructureDefinition/comment}`n"
$ws.Range("L20").Value = "Comment"
$ws.Range("M20").Value = "Used to describe the social background of a patient affiliated to CARA"
$ws.Range("N20").ClearContents()
$ws.Range("O20").ClearContents()
$ws.Range("P20").Value = ""
$ws.Range("Q20").ClearContents()
$ws.Range("R20").Value = ""
$ws.Range("S20").Value = ""
$ws.Range("T20").Value = ""
$ws.Range("U20").Value = ""
$ws.Range("V20").Value = ""
$ws.Range("W20").Value = ""
$ws.Range("X20").Value = ""
$ws.Range("Y20").Value = ""
$ws.Range("Z20").Value = ""
$ws.Range("AA20").Value = ""
$ws.Range("AB20").Value = ""
$ws.Range("AC20").Value = ""
$ws.Range("AD20").Value = ""
$ws.Range("AE20").Value = ""
$ws.Range("AF20").Value = "Element.extension"
$ws.Range("AG20").Value = "0"
$ws.Range("AH20").Value = "*"
$ws.Range("AI20").Value = "ele-1`n"
$ws.Range("AJ20").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$ws.Range("AK20").Value = ""
$ws.Range("AL20").Value = ""
$ws.Range("AM20").Value = ""

# ---------------------------------------------------------------------
# Column sizing tweaks that result from the new/longer content
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 38.3828125
$ws.Columns.Item(3).ColumnWidth = 9.37890625
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(28).ColumnWidth = 9.2890625
$ws.Columns.Item(31).ColumnWidth = 5.484375
